$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.649191975593567
$ws.Range("B1").Value = 2.517845869064331
$ws.Range("C1").Value = 2.709914922714233
$ws.Range("D1").Value = 3.045075654983521
$ws.Range("E1").Value = 3.545687198638916
